$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "End date" estimates in G2:G5 - times will be decided later
$ws.Range("G2:G5").ClearContents()

# Update the active selection to B9
$ws.Range("B9").Select()
